$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (modelo, politica) before the existing "full" column (C)
# Existing layout: A=nome, B=preco, C=full, D=tipo, E=link
# Target layout:   A=nome, B=preco, C=modelo, D=politica, E=full, F=tipo, G=link
$ws.Columns("C:D").Insert()

# Header row
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Row 2 - Fonte Carregador Jfa Bob Storm 90a
$ws.Range("C2").Value = "FONTE 90 BOB"
$ws.Range("D2").Value = "Igual"
$ws.Range("F2").Value = "classico"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:1165626720#searchVariation=MLB21562641&position=44&search_layout=grid&type=product&tracking_id=a8149985-e671-49c4-9f66-79445fc51396"

# Row 3 - Controle Longa Distancia Jfa Acqua 1200 (MLB27685629)
$ws.Range("C3").Value = "Sem Modelo"
$ws.Range("D3").Value = ""
$ws.Range("F3").Value = "classico"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27685629?pdp_filters=seller_id:1165626720#searchVariation=MLB27685629&position=34&search_layout=grid&type=product&tracking_id=75974dfc-751d-4169-a163-00ed9f3cf548"

# Row 4 - Controle Remoto Universal Longa Distancia Jfa K1200 Azul
$ws.Range("C4").Value = "Sem Modelo"
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-azul/p/MLB28722231?pdp_filters=seller_id:1165626720#searchVariation=MLB28722231&position=24&search_layout=grid&type=product&tracking_id=ca8e3fdd-4279-4982-948e-579f8729449e"

# Row 5 - Fonte Automotiva 40 Amperes
$ws.Range("C5").Value = "FONTE 40A"
$ws.Range("D5").Value = "Igual"
$ws.Range("F5").Value = "classico"
$ws.Range("G5").Value = "https://www.mercadolivre.com.br/fonte-automotiva-40-amperes-jfa-storm-red-line-cca-sci-smart-cor-preto/p/MLB21621306?pdp_filters=seller_id:1165626720#searchVariation=MLB21621306&position=1&search_layout=grid&type=product&tracking_id=5126fbcf-05ab-4228-8346-ef45f807253f"

# Row 6 - Fonte Carregador De Bateria Jfa 70a Lite Slim Bivolt
$ws.Range("C6").Value = "FONTE 70A LITE"
$ws.Range("D6").Value = "Acima"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "classico"
$ws.Range("G6").Value = "https://produto.mercadolivre.com.br/MLB-3715575332-fonte-carregador-de-bateria-jfa-70a-lite-slim-bivolt-_JM#position%3D35%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Dc0f7b93b-50af-43cd-9269-a325b03cbe9c"

# Row 7 - Controle Longa Distancia Jfa Acqua 1200 (MLB27687422)
$ws.Range("C7").Value = "Sem Modelo"
$ws.Range("D7").Value = ""
$ws.Range("F7").Value = "classico"
$ws.Range("G7").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:1165626720#searchVariation=MLB27687422&position=1&search_layout=grid&type=product&tracking_id=0a55ba12-2aaf-47b8-8be7-eb13713785bb"
